# Add team record (Wins/Losses/Ties) columns to the NYY 1993 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): AD1=Wins, AE1=Losses, AF1=Ties, matching the style
# of the other header cells (bold/centered/bordered -> same as A1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$hdr = $ws.Range("AD1:AF1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108  # xlCenter
$hdr.VerticalAlignment = -4160    # xlTop
$hdr.Borders.LineStyle = 1        # xlContinuous (thin box border)

# Data rows 2-43: team record for every player row (NYY finished 1993
# with an 88-74-0 record).
$lastRow = 43
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 88  # column AD
    $ws.Cells.Item($r, 31).Value = 74  # column AE
    $ws.Cells.Item($r, 32).Value = 0   # column AF
}
